$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SKYS")

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D:D").Insert()

# Copy number formats from the (now-shifted) old column into the new column D
# so the new cells retain the same style as before (date / number formats)
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate new column D ("latest quarter") values
$ws.Range("D7").Value = 43281
$ws.Range("D8").Value = 33200
$ws.Range("D9").Value = 15000
$ws.Range("D10").Value = 18200
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 8800
$ws.Range("D18").Value = 24400
$ws.Range("D20").Value = 2100
$ws.Range("D21").Value = "NA"
$ws.Range("D22").Value = 8500
$ws.Range("D23").Value = 17900
$ws.Range("D24").Value = 9700
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 8200
$ws.Range("D27").Value = 8200
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -2100
$ws.Range("D33").Value = 8200
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 8200
$ws.Range("D38").Value = 43281
$ws.Range("D41").Value = 68100
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 46700
$ws.Range("D44").Value = 300
$ws.Range("D45").Value = 40500
$ws.Range("D46").Value = 155600
$ws.Range("D47").Value = 5200
$ws.Range("D48").Value = 403300
$ws.Range("D49").Value = "NA"
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 48100
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 612200
$ws.Range("D57").Value = 21800
$ws.Range("D58").Value = 24000
$ws.Range("D59").Value = 132400
$ws.Range("D60").Value = 178200
$ws.Range("D61").Value = 248700
$ws.Range("D62").Value = 73100
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 505100
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 107100
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 107100
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43281
$ws.Range("D81").Value = 8200
$ws.Range("D83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 0
$ws.Range("D91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 0
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 0
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 0
